$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '329.87'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '6.77%'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '40.12'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '6.90%'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.267'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.83%'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08100'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2.91%'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.522'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '2.22%'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.653'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '4.92%'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.931'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '1.10%'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9361'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-0.07%'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1351'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '21.33%'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1969'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.09%'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09087'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.71%'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03503'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '6.06%'

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.12%'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001409'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1.64%'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.006159'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.82%'

$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.383'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-6.40%'

$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3519'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '3.15%'

$ws.Range("B20").Value = 'MCDex'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.486'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.53%'

$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1314'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '2.52%'

$ws.Range("B22").Value = 'ZBToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2572'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '2.09%'

$ws.Range("B23").Value = 'CoinExToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04448'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.87%'

$ws.Range("B24").Value = 'BitKan'
$ws.Range("C24").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001224'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.92%'

$ws.Range("B25").Value = 'HotbitToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004320'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-5.59%'

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-5.21%'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0003994'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '0.04%'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02507'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '13.68%'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05199'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '2.13%'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007712'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '3.37%'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1428'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '5.49%'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.009191'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '4.87%'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002163'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '1.33%'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.009010'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '4.42%'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006629'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '1.23%'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.06%'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003345'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '16.84%'

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '147.76%'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002103'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.06%'

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.06%'
